# Insert a new row before row 234; this shifts existing rows 234-318 down to 235-319.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(234).Insert()

# Populate the newly inserted row 234 with the new data record.
$ws.Range("A234").Value = 9
$ws.Range("B234").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C234").Value = "Metropolitana"
$ws.Range("D234").Value = 44559
$ws.Range("E234").Value = 13
$ws.Range("F234").Value = 100112039
$ws.Range("G234").Value = "Ciboulette"
$ws.Range("H234").Value = "Sin especificar"
$ws.Range("I234").Value = "Primera"
$ws.Range("J234").Value = 160
$ws.Range("K234").Value = 1000
$ws.Range("L234").Value = 1200
$ws.Range("M234").Value = 1100
$ws.Range("N234").Value = "`$/docena de atados"
$ws.Range("O234").Value = "Región Metropolitana"
$ws.Range("P234").Value = 367
$ws.Range("Q234").Value = 3
$ws.Range("R234").Value = "Hortaliza"
